$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently holds years 2000..2020 in rows 2..22.
# Target state: years 2010..2022 in rows 2..14.
# Step 1: drop the obsolete 2000..2009 rows (old rows 2..11) -- this shifts
# everything below up by 10 rows, so old row 12 (2010) becomes new row 2,
# ... old row 22 (2020) becomes new row 12 (preserving its cell records,
# including the two literal empty-string cells in B/C exactly as they were).
$ws.Range("A2:F11").EntireRow.Delete()

# Step 2: append the two new years, 2021 and 2022, as rows 13 and 14.
$ws.Range("A13").Value = "2021年"
$ws.Range("D13").Value = 333.2373
$ws.Range("E13").Value = 117.6526
$ws.Range("F13").Value = 77.2761

$ws.Range("A14").Value = "2022年"
$ws.Range("D14").Value = 365.3613
$ws.Range("E14").Value = 124.2479
$ws.Range("F14").Value = 86.2165

# Match the year-label column styling (bold, centered, bordered) used by
# every other cell in column A -- copy the format from the row right above
# (2020, row 12) so the new rows reuse the same cell style instead of
# minting a near-duplicate one.
$ws.Range("A12").Copy()
$ws.Range("A13:A14").PasteSpecial(-4122)
